# Generate Report for Handoff
#
# The "7d53c785-489c-4cfd-b166-a85032a26e52" entry has finished translation
# and is now ready to be handed off again, while the
# "9da099d2-f805-4acf-957c-9dc68465a39b" entry is (still) in translation.
# The report rows are re-sorted so the "in translation" item now comes
# first (row 2) and the "ready for handoff" item comes second (row 3);
# the handoff status/time-stamps for the 7d53c785 item are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "9da099d2-f805-4acf-957c-9dc68465a39b.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"

$ov.Range("A3").Value = "7d53c785-489c-4cfd-b166-a85032a26e52.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/7d53c785-489c-4cfd-b166-a85032a26e52.md", "", "", "9da099d2-f805-4acf-957c-9dc68465a39b.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/9da099d2-f805-4acf-957c-9dc68465a39b.md", "", "", "7d53c785-489c-4cfd-b166-a85032a26e52.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "9da099d2-f805-4acf-957c-9dc68465a39b.md"
$zh.Range("B2").Value = "In Translation"
$zh.Range("C2").Value = "9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-09 10:21:15"
$zh.Range("G2").Value = "0001-01-01 00:00:00"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "7d53c785-489c-4cfd-b166-a85032a26e52.md"
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-09 10:22:05"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/7d53c785-489c-4cfd-b166-a85032a26e52.md", "", "", "9da099d2-f805-4acf-957c-9dc68465a39b.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5eea6b846229514ed432272182590afc00322332/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.zh-cn.xlf", "", "", "9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/9da099d2-f805-4acf-957c-9dc68465a39b.md", "", "", "7d53c785-489c-4cfd-b166-a85032a26e52.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5eea6b846229514ed432272182590afc00322332/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.zh-cn.xlf", "", "", "7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "9da099d2-f805-4acf-957c-9dc68465a39b.md"
$de.Range("B2").Value = "In Translation"
$de.Range("C2").Value = "9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.de-de.xlf"
$de.Range("D2").Value = "2016-03-09 10:21:22"
$de.Range("G2").Value = "0001-01-01 00:00:00"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "7d53c785-489c-4cfd-b166-a85032a26e52.md"
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.de-de.xlf"
$de.Range("D3").Value = "2016-03-09 10:22:09"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/7d53c785-489c-4cfd-b166-a85032a26e52.md", "", "", "9da099d2-f805-4acf-957c-9dc68465a39b.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e15f4fcfefa9d9453f8fa16ac43875732a897bb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.de-de.xlf", "", "", "9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/e2e/9da099d2-f805-4acf-957c-9dc68465a39b.md", "", "", "7d53c785-489c-4cfd-b166-a85032a26e52.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e15f4fcfefa9d9453f8fa16ac43875732a897bb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9da099d2-f805-4acf-957c-9dc68465a39b.171e89cf0b196b8ed524981b4edcb94610c2bbf2.de-de.xlf", "", "", "7d53c785-489c-4cfd-b166-a85032a26e52.cb082305a6b3eb9a033914ec5e6e841222ee657c.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ce8d5ebbd03dc225c450edbcdc030013dc1cb1c3/.localization-config", "", "", ".localization-config")
